# TICC-123: Deprecated Peppol AS2 v1
#
# On the "Transport Profile" sheet, row 3 is the AS2 1.0 profile
# (busdox-transport-as2-ver1p0). Mark it as deprecated, the same way the
# already-deprecated AS4 v1 row (row 4) is represented: a TRUE() formula in
# the "Deprecated" column, and fill in the "Deprecated since?" column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Formula = "=TRUE"
$ws.Range("F3").Value = 7

# Editor cursor/selection also moved from C7 to A7 in the authored file.
$ws.Range("A7").Select()
